# Update the Implementation Guide metadata and element sheets.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet -------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# URL: moved from path-based to subdomain-based hostname
$wsMeta.Range("B2").Value = "https://molic-avc.gabriellesantosleandro.com/StructureDefinition/live-alone-pre-stroke-extension"

# Date: updated publication timestamp
$wsMeta.Range("B8").Value = "2023-08-16T00:27:03-03:00"

# --- "Elements" sheet ---------------------------------------------------
$wsElements = $wb.Worksheets.Item("Elements")

# Extension.url fixed value mirrors the canonical URL above
$wsElements.Range("R4").Value = "https://molic-avc.gabriellesantosleandro.com/StructureDefinition/live-alone-pre-stroke-extension"

# ValueSet URL: same hostname restructuring as above
$wsElements.Range("Z5").Value = "https://molic-avc.gabriellesantosleandro.com/ValueSet/live-alone-pre-stroke"

# Narrow column Z (26) to fit the shorter URL text (target stored width
# 70.19921875 chars; this runtime quantizes ColumnWidth writes to 1/6-char
# steps, so 69.33333333333333 is the closest input that lands on it)
$wsElements.Columns.Item(26).ColumnWidth = 69.33333333333333
